$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "298.08"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "1.35%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "42.11"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "3.98%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.011"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-0.28%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07520"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "2.64%"
$ws.Range("B6").Value = "FTXToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.593"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "1.74%"
$ws.Range("B7").Value = "MXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9297"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-0.17%"
$ws.Range("B8").Value = "BTSEToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "2.401"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "2.01%"
$ws.Range("B9").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C9").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1190"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "2.52%"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1833"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "4.18%"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08986"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "2.69%"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.04134"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-4.96%"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.1050"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.53%"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001280"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "2.71%"
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.005794"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-1.51%"
$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.337"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-0.33%"
$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.366"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "1.82%"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "0.99%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "8.306"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "5.61%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1390"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "0.03%"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "12.04%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.04098"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "4.72%"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "0.47%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.003894"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "6.28%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0001301"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "8.56%"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02398"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "3.11%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05227"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "2.31%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.006774"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "17.38%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007815"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-0.51%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1327"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "2.52%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.007405"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "0.55%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.007113"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-1.81%"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "1.92%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006560"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "6.15%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000750"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.20%"
$ws.Range("B48").Value = "BOLO"
$ws.Range("C48").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.04904"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "0.24%"
$ws.Range("B49").Value = "CoinbaseStockToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.004202"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "0.00%"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.20%"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.20%"
